$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings (e.g. "590.03")
# must be forced to Text format first, otherwise Excel auto-converts them
# to actual numbers (losing the intended text representation / trailing zeros).
$textCells = @(
    "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D17", "D21", "D22", "D23", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D35", "D37", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49"
)
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Row 2
$ws.Range('D2').Value = '70.501.30'
$ws.Range('E2').Value = '  -2.89%  '

# Row 3
$ws.Range('D3').Value = '3.851.62'
$ws.Range('E3').Value = '  -3.06%  '

# Row 4
$ws.Range('E4').Value = '  +0.18%  '

# Row 5
$ws.Range('D5').Value = '590.03'
$ws.Range('E5').Value = '  +0.84%  '

# Row 6
$ws.Range('D6').Value = '166.01'
$ws.Range('E6').Value = '  +5.25%  '

# Row 7
$ws.Range('D7').Value = '0.666'
$ws.Range('E7').Value = '  -1.88%  '

# Row 8
$ws.Range('E8').Value = '  +0.28%  '

# Row 9
$ws.Range('D9').Value = '0.744'
$ws.Range('E9').Value = '  -0.57%  '

# Row 10
$ws.Range('D10').Value = '0.173'
$ws.Range('E10').Value = '  +3.32%  '

# Row 11
$ws.Range('D11').Value = '52.99'
$ws.Range('E11').Value = '  -0.60%  '

# Row 12
$ws.Range('D12').Value = '0.0000319'
$ws.Range('E12').Value = '  +0.16%  '

# Row 13
$ws.Range('D13').Value = '11.24'
$ws.Range('E13').Value = '  +3.38%  '

# Row 14
$ws.Range('D14').Value = '4.474.90'
$ws.Range('E14').Value = '  -2.82%  '

# Row 15
$ws.Range('D15').Value = '3.872.41'
$ws.Range('E15').Value = '  -2.45%  '

# Row 16
$ws.Range('E16').Value = '  +0.99%  '

# Row 17
$ws.Range('D17').Value = '13.74'
$ws.Range('E17').Value = '  -2.11%  '

# Row 18
$ws.Range('E18').Value = '  -6.42%  '

# Row 19
$ws.Range('E19').Value = '  -2.22%  '

# Row 20
$ws.Range('D20').Value = '70.421.57'
$ws.Range('E20').Value = '  -2.55%  '

# Row 21
$ws.Range('D21').Value = '433.96'
$ws.Range('E21').Value = '  +0.08%  '

# Row 22
$ws.Range('D22').Value = '4.70'
$ws.Range('E22').Value = '  +0.49%  '

# Row 23
$ws.Range('D23').Value = '93.94'
$ws.Range('E23').Value = '  -2.08%  '

# Row 24
$ws.Range('E24').Value = '  -4.97%  '

# Row 25
$ws.Range('E25').Value = '  -4.29%  '

# Row 26
$ws.Range('D26').Value = '4.06'
$ws.Range('E26').Value = '  -9.12%  '

# Row 27
$ws.Range('D27').Value = '10.92'
$ws.Range('E27').Value = '  -1.74%  '

# Row 28
$ws.Range('E28').Value = '  +0.10%  '

# Row 29
$ws.Range('D29').Value = '10.22'
$ws.Range('E29').Value = '  -4.53%  '

# Row 30
$ws.Range('D30').Value = '34.92'
$ws.Range('E30').Value = '  -4.20%  '

# Row 31
$ws.Range('D31').Value = '7.92'
$ws.Range('E31').Value = '  +1.21%  '

# Row 32
$ws.Range('D32').Value = '13.41'
$ws.Range('E32').Value = '  -1.24%  '

# Row 33
$ws.Range('D33').Value = '48.37'
$ws.Range('E33').Value = '  -0.95%  '

# Row 34
$ws.Range('E34').Value = '  -4.93%  '

# Row 35
$ws.Range('D35').Value = '69.20'
$ws.Range('E35').Value = '  +0.41%  '

# Row 36
$ws.Range('D36').Value = '0.0₃0964'
$ws.Range('E36').Value = '  +11.88%  '

# Row 37
$ws.Range('D37').Value = '614.55'
$ws.Range('E37').Value = '  -9.81%  '

# Row 38
$ws.Range('D38').Value = '0.418'
$ws.Range('E38').Value = '  -4.11%  '

# Row 39
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  +0.02%  '

# Row 40
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.06%  '

# Row 41
$ws.Range('E41').Value = '  -2.44%  '

# Row 42
$ws.Range('B42').Value = 'ThetaToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D42').Value = '3.24'
$ws.Range('E42').Value = '  -3.40%  '

# Row 43
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').Value = '3.20'
$ws.Range('E43').Value = '  +25.62%  '

# Row 44
$ws.Range('D44').Value = '0.0464'
$ws.Range('E44').Value = '  -4.20%  '

# Row 45
$ws.Range('D45').Value = '9.98'
$ws.Range('E45').Value = '  -7.41%  '

# Row 46
$ws.Range('D46').Value = '2.67'
$ws.Range('E46').Value = '  +0.65%  '

# Row 47
$ws.Range('D47').Value = '0.142'
$ws.Range('E47').Value = '  -4.07%  '

# Row 48
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').Value = '3.29'
$ws.Range('E48').Value = '  -2.45%  '

# Row 49
$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').Value = '2.78'
$ws.Range('E49').Value = '  -17.09%  '

# Row 50
$ws.Range('D50').Value = '2.829.72'
$ws.Range('E50').Value = '  +2.00%  '

# Row 51
$ws.Range('E51').Value = '  +1.13%  '

# Restore General number format on the cells we temporarily forced to Text
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "General"
}